$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Section: "Przelewy z mbanku" (rows 12-18)
# ---------------------------------------------------------------------
$ws.Range("A12").Value = "Przelewy z mbanku"

$ws.Range("A13").Value = "data"
$ws.Range("B13").Value = "tytuł"
$ws.Range("C13").Value = "kwota"
$ws.Range("A13:C13").Font.Bold = $true
$ws.Range("A13:C13").Borders.LineStyle = 1
$ws.Range("C13").NumberFormat = "0.00"

$ws.Range("A14").Value = "26-07-2012"
$ws.Range("B14").Value = "TYSIAK"
$ws.Range("C14").Value = 1000

$ws.Range("A15").Value = "27-08-2012"
$ws.Range("B15").Value = 2500
$ws.Range("C15").Value = 2500

$ws.Range("A16").Value = "19-09-2012"
$ws.Range("B16").Value = "tysiąć"
$ws.Range("C16").Value = 1000

$ws.Range("A17").Value = "04-10-2012"
$ws.Range("B17").Value = "damy radę"
$ws.Range("C17").Value = 1000

$ws.Range("A18").Value = "23-11-2012"
$ws.Range("B18").Value = 500
$ws.Range("C18").Value = 500

$ws.Range("A14:C18").Borders.LineStyle = 1
$ws.Range("B14:B17").HorizontalAlignment = -4131
$ws.Range("B18").HorizontalAlignment = -4131
$ws.Range("B18").NumberFormat = '#,##0\ "zł";[Red]\-#,##0\ "zł"'

# ---------------------------------------------------------------------
# Section: "Przelewy z polbanku" (rows 20-23)
# ---------------------------------------------------------------------
$ws.Range("A20").Value = "Przelewy z polbanku"

$ws.Range("A21").Value = "data"
$ws.Range("B21").Value = "tytuł"
$ws.Range("C21").Value = "kwota"
$ws.Range("A21:C21").Font.Bold = $true
$ws.Range("A21:C21").Borders.LineStyle = 1
$ws.Range("C21").NumberFormat = "0.00"

$ws.Range("A22").Value = "??"
$ws.Range("B22").Value = "??"
$ws.Range("C22").Value = "??"
$ws.Range("A23").Value = "??"
$ws.Range("B23").Value = "??"
$ws.Range("C23").Value = "??"

$ws.Range("A22:C23").Borders.LineStyle = 1
$ws.Range("B22:B23").HorizontalAlignment = -4131

# ---------------------------------------------------------------------
# Section: "Przelewy z eb" (rows 25-27)
# ---------------------------------------------------------------------
$ws.Range("A25").Value = "Przelewy z eb"

$ws.Range("A26").Value = "data"
$ws.Range("B26").Value = "tytuł"
$ws.Range("C26").Value = "kwota"
$ws.Range("A26:C26").Font.Bold = $true
$ws.Range("A26:C26").Borders.LineStyle = 1
$ws.Range("C26").NumberFormat = "0.00"

$ws.Range("A27").Value = "15-10-2012"
$ws.Range("B27").Value = "tysiąc"
$ws.Range("C27").Value = 1000
$ws.Range("A27:B27").Borders.LineStyle = 1
$ws.Range("C27").Borders.LineStyle = 1
$ws.Range("C27").NumberFormat = "0.00"

# ---------------------------------------------------------------------
# New trailing blank row (row 31) mirrors the rest of the table
# ---------------------------------------------------------------------
$ws.Range("C31").Value = $ws.Range("C30").Value
$ws.Range("C31").NumberFormat = $ws.Range("C30").NumberFormat
$ws.Range("D31").Value = $ws.Range("D30").Value
$ws.Range("D31").NumberFormat = $ws.Range("D30").NumberFormat
$ws.Range("E31").Value = $ws.Range("E30").Value
$ws.Range("E31").NumberFormat = $ws.Range("E30").NumberFormat
$ws.Range("C31").ClearContents()
$ws.Range("D31").ClearContents()
$ws.Range("E31").ClearContents()

# ---------------------------------------------------------------------
# Selection / view tweaks
# ---------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("I26").Select()
